$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 3 was "ammo_magnum_300"-style "AP" ammo; it's now re-purposed as the
# Deer Hunter load (.338 Federal), classified as "DMG" ammo instead of "AP".
$ws.Range("C3").Value = 4000
$ws.Range("B3").Value = "DMG"
$ws.Range("H3").Value = 10.7
$ws.Range("J3").Value = 9.0949999999999989

# Highlight the changed row with the plain Accent2 theme color (no tint),
# used to flag freshly-edited entries.
$ws.Range("A3:B3").Font.ThemeColor = 6
$ws.Range("E3").Font.ThemeColor = 6

$ws.Range("E11").Select() | Out-Null
